$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.208684666666667
$ws.Range("H2").Value = 3.626054
$ws.Range("I2").Value = 0.01462795763842055
$ws.Range("J2").Value = 0.01462795763842055
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03257366666666667
$ws.Range("N2").Value = 0.097721
$ws.Range("O2").Value = 0.001227793554179957
$ws.Range("P2").Value = 0.001227793554179957
$ws.Range("Q2").Value = 0.03937129143711111
$ws.Range("R2").Value = 0.354341622934
$ws.Range("S2").Value = 0.00001796011209927023
$ws.Range("T2").Value = 0.00001796011209927023
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.208684666666667
$ws.Range("H3").Value = 3.626054
$ws.Range("I3").Value = 0.01462795763842055
$ws.Range("J3").Value = 0.01462795763842055
$ws.Range("O3").Value = 0.7662385783512358
$ws.Range("P3").Value = 0.7662385783512359
$ws.Range("Q3").Value = 24.57074503764866
$ws.Range("R3").Value = 221.136705338838
$ws.Range("S3").Value = 0.01120850546504547
$ws.Range("T3").Value = 0.01120850546504547
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.208684666666667
$ws.Range("H4").Value = 3.626054
$ws.Range("I4").Value = 0.01462795763842055
$ws.Range("J4").Value = 0.01462795763842055
$ws.Range("M4").Value = 6.169174666666667
$ws.Range("N4").Value = 18.507524
$ws.Range("O4").Value = 0.2325336280945842
$ws.Range("P4").Value = 0.2325336280945842
$ws.Range("Q4").Value = 7.456586825588444
$ws.Range("R4").Value = 67.109281430296
$ws.Range("S4").Value = 0.003401492061275817
$ws.Range("T4").Value = 0.003401492061275817
$ws.Range("I5").Value = 0.0626664797952065
$ws.Range("J5").Value = 0.06266647979520648
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03257366666666667
$ws.Range("N5").Value = 0.097721
$ws.Range("O5").Value = 0.001227793554179957
$ws.Range("P5").Value = 0.001227793554179957
$ws.Range("Q5").Value = 0.1686674449257778
$ws.Range("R5").Value = 1.518007004332
$ws.Range("S5").Value = 0.00007694149995570307
$ws.Range("T5").Value = 0.00007694149995570305
$ws.Range("I6").Value = 0.0626664797952065
$ws.Range("J6").Value = 0.06266647979520648
$ws.Range("O6").Value = 0.7662385783512358
$ws.Range("P6").Value = 0.7662385783512359
$ws.Range("S6").Value = 0.04801747438855546
$ws.Range("T6").Value = 0.04801747438855546
$ws.Range("I7").Value = 0.0626664797952065
$ws.Range("J7").Value = 0.06266647979520648
$ws.Range("M7").Value = 6.169174666666667
$ws.Range("N7").Value = 18.507524
$ws.Range("O7").Value = 0.2325336280945842
$ws.Range("P7").Value = 0.2325336280945842
$ws.Range("Q7").Value = 31.94417561202312
$ws.Range("R7").Value = 287.497580508208
$ws.Range("S7").Value = 0.01457206390669532
$ws.Range("T7").Value = 0.01457206390669532
$ws.Range("G8").Value = 76.16218566666667
$ws.Range("H8").Value = 228.486557
$ws.Range("I8").Value = 0.9217434921665711
$ws.Range("J8").Value = 0.921743492166571
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03257366666666667
$ws.Range("N8").Value = 0.097721
$ws.Range("O8").Value = 0.001227793554179957
$ws.Range("P8").Value = 0.001227793554179957
$ws.Range("Q8").Value = 2.480881648510778
$ws.Range("R8").Value = 22.327934836597
$ws.Range("S8").Value = 0.00113171071828944
$ws.Range("T8").Value = 0.00113171071828944
$ws.Range("G9").Value = 76.16218566666667
$ws.Range("H9").Value = 228.486557
$ws.Range("I9").Value = 0.9217434921665711
$ws.Range("J9").Value = 0.921743492166571
$ws.Range("O9").Value = 0.7662385783512358
$ws.Range("P9").Value = 0.7662385783512359
$ws.Range("Q9").Value = 1548.262915162648
$ws.Range("R9").Value = 13934.36623646383
$ws.Range("S9").Value = 0.7062754230422169
$ws.Range("T9").Value = 0.7062754230422169
$ws.Range("G10").Value = 76.16218566666667
$ws.Range("H10").Value = 228.486557
$ws.Range("I10").Value = 0.9217434921665711
$ws.Range("J10").Value = 0.921743492166571
$ws.Range("M10").Value = 6.169174666666667
$ws.Range("N10").Value = 18.507524
$ws.Range("O10").Value = 0.2325336280945842
$ws.Range("P10").Value = 0.2325336280945842
$ws.Range("Q10").Value = 469.8578263727632
$ws.Range("R10").Value = 4228.720437354868
$ws.Range("S10").Value = 0.2143363584060647
$ws.Range("T10").Value = 0.2143363584060647
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.07949433333333333
$ws.Range("H11").Value = 0.238483
$ws.Range("I11").Value = 0.0009620703998019471
$ws.Range("J11").Value = 0.000962070399801947
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.03257366666666667
$ws.Range("N11").Value = 0.097721
$ws.Range("O11").Value = 0.001227793554179957
$ws.Range("P11").Value = 0.001227793554179957
$ws.Range("Q11").Value = 0.002589421915888889
$ws.Range("R11").Value = 0.023304797243
$ws.Range("S11").Value = 0.000001181223835544165
$ws.Range("T11").Value = 0.000001181223835544165
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.07949433333333333
$ws.Range("H12").Value = 0.238483
$ws.Range("I12").Value = 0.0009620703998019471
$ws.Range("J12").Value = 0.000962070399801947
$ws.Range("O12").Value = 0.7662385783512358
$ws.Range("P12").Value = 0.7662385783512359
$ws.Range("Q12").Value = 1.616000475672333
$ws.Range("R12").Value = 14.544004281051
$ws.Range("S12").Value = 0.0007371754554180489
$ws.Range("T12").Value = 0.0007371754554180489
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.07949433333333333
$ws.Range("H13").Value = 0.238483
$ws.Range("I13").Value = 0.0009620703998019471
$ws.Range("J13").Value = 0.000962070399801947
$ws.Range("M13").Value = 6.169174666666667
$ws.Range("N13").Value = 18.507524
$ws.Range("O13").Value = 0.2325336280945842
$ws.Range("P13").Value = 0.2325336280945842
$ws.Range("Q13").Value = 0.4904144273435556
$ws.Range("R13").Value = 4.413729846092
$ws.Range("S13").Value = 0.0002237137205483539
$ws.Range("T13").Value = 0.0002237137205483539
